$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A55").Value = 20220201
$ws.Range("B55").Value = 2228.2564673073002
$ws.Range("C55").Value = 2224.4699999999998
$ws.Range("D55").Formula = "=100*(B55-C55)/C55"
$ws.Range("E55").Value = 180
$ws.Range("F55").Value = "CRM OPENED 20220118"
